# Fruta / hortaliza, semanal
# Insert this week's new price records (date 44627) at the top of the
# Membrillo price history table (row 50), pushing the existing history
# down by two rows (old row 50 -> new row 52, etc.) and growing the used
# range from A1:T67 to A1:T69.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 50..67 down by two rows to make room for the new week's data.
$ws.Rows("50:51").Insert()

# New row 50: Especial
$ws.Range("A50").Value = 3
$ws.Range("B50").Value = "Femacal de La Calera"
$ws.Range("C50").Value = "Coquimbo"
$ws.Range("D50").Value = 44627
$ws.Range("E50").Value = 5
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100104
$ws.Range("H50").Value = "Frutos de pepita"
$ws.Range("I50").Value = 100104003
$ws.Range("J50").Value = "Membrillo"
$ws.Range("K50").Value = "Champion"
$ws.Range("L50").Value = "Especial"
$ws.Range("M50").Value = 65
$ws.Range("N50").Value = 15000
$ws.Range("O50").Value = 15000
$ws.Range("P50").Value = 15000
$ws.Range("Q50").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 833
$ws.Range("T50").Value = 18

# New row 51: Primera
$ws.Range("A51").Value = 3
$ws.Range("B51").Value = "Femacal de La Calera"
$ws.Range("C51").Value = "Coquimbo"
$ws.Range("D51").Value = 44627
$ws.Range("E51").Value = 5
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100104
$ws.Range("H51").Value = "Frutos de pepita"
$ws.Range("I51").Value = 100104003
$ws.Range("J51").Value = "Membrillo"
$ws.Range("K51").Value = "Champion"
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 60
$ws.Range("N51").Value = 13000
$ws.Range("O51").Value = 13000
$ws.Range("P51").Value = 13000
$ws.Range("Q51").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R51").Value = "Región de O'Higgins"
$ws.Range("S51").Value = 722
$ws.Range("T51").Value = 18
